$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet1: make room for a new "Sentences removed" column before the old
# "Reviews removed" column (D). Only rows 1 and (eventually) 3 get touched;
# row 2's new D cell is cleared right away so it stays completely blank,
# matching the target (no stray <c r="D2"/> placeholder).
$ws.Range("D1:D2").Insert(-4161)
$ws.Range("D2").Clear()

# --- Row 3: brand new "Per review" stats row. Written before the row-1/row-2
# labels below so its new shared string ("Per review") gets the lowest new
# shared-string index, matching the target workbook's string order.
$ws.Range("A3").Value = "Per review"
$ws.Range("B3").Formula = "=1/3"
$ws.Range("B3").NumberFormat = "0.00"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 33808
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").VerticalAlignment = -4108
$ws.Range("D3").WrapText = $true
$ws.Range("E3").Value = 1046
$ws.Range("F3").Value = 14400
$ws.Range("G3").Formula = "=F3/(E3+F3)"
$ws.Range("G3").NumberFormat = "0.00"

# --- Row 1: new "Remove from" header in column A
$ws.Range("A1").Value = "Remove from"

# --- Row 2: new "Corpus" label in column A
$ws.Range("A2").Value = "Corpus"

# --- Row 1: new "Sentences removed" header in column D
$ws.Range("D1").Value = "Sentences removed"

# Row heights for the two data rows grew slightly (16 -> 17) after the edits.
$ws.Rows("2:3").RowHeight = 17

# --- Sheet view / selection bookkeeping: Sheet1 becomes the active tab
# (previously "Top subj" / sheet index 3 was active), and the cell cursor on
# Sheet1 ends up at H7.
$ws.Activate()
$ws.Range("H7").Select()
